# fix(docx): fix OOXMLValidator error on KeywordTok output
#
# wml.xsd's CT_RPr sequence requires run-property toggle elements such as
# <w:b/> / <w:i/> to come before <w:color/>. Several of the Pandoc
# "*Tok" character styles in styles.xml had <w:color/> emitted before
# <w:b/>/<w:i/>, which OOXMLValidatorCLI flags as
# Sch_UnexpectedElementContentExpectingComplex even though xmllint stays
# quiet about it. Re-assert the bold/italic toggles on each affected
# style so the run properties get re-serialized in schema order
# (toggle elements before <w:color/>), without changing any of the
# actual formatting (same color, same bold/italic state).

$d = $word.ActiveDocument

# styleId -> whether Bold / Italic should be (re-)applied, in this order.
$fixes = @(
    @{ Name = "KeywordTok";       Bold = $true;  Italic = $false },
    @{ Name = "ImportTok";        Bold = $true;  Italic = $false },
    @{ Name = "CommentTok";       Bold = $false; Italic = $true  },
    @{ Name = "DocumentationTok"; Bold = $false; Italic = $true  },
    @{ Name = "AnnotationTok";    Bold = $true;  Italic = $true  },
    @{ Name = "CommentVarTok";    Bold = $true;  Italic = $true  },
    @{ Name = "ControlFlowTok";   Bold = $true;  Italic = $false },
    @{ Name = "InformationTok";   Bold = $true;  Italic = $true  },
    @{ Name = "WarningTok";       Bold = $true;  Italic = $true  },
    @{ Name = "AlertTok";         Bold = $true;  Italic = $false },
    @{ Name = "ErrorTok";         Bold = $true;  Italic = $false }
)

foreach ($fix in $fixes) {
    $style = $d.Styles($fix.Name)
    if ($fix.Bold) {
        $style.Font.Bold = $true
    }
    if ($fix.Italic) {
        $style.Font.Italic = $true
    }
}
